# Generate Report for Handback
# Update status/handback datetime/error detail for zh-cn and de-de target
# locales, and mirror the status on the Overview sheet. Also widen a couple
# of columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: mirrors the zh-cn / de-de status in columns E (zh-cn) and F (de-de) ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$overview.Columns.Item(5).ColumnWidth = 29.1666666666667
$overview.Columns.Item(6).ColumnWidth = 29.1666666666667

# --- zh-cn sheet ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$zhcn.Range("K2").Value = "2016-10-19 12:03:54"
$zhcn.Range("K3").Value = "2016-10-19 12:03:54"

$zhcn.Range("P2").Value = ""

$zhcn.Columns.Item(3).ColumnWidth = 29.1666666666667
$zhcn.Columns.Item(16).ColumnWidth = 12.8333333333333

# --- de-de sheet ---
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Range("K2").Value = "2016-10-19 12:04:13"
$dede.Range("K3").Value = "2016-10-19 12:04:13"

$dede.Range("P2").Value = ""

$dede.Columns.Item(3).ColumnWidth = 29.1666666666667
$dede.Columns.Item(16).ColumnWidth = 12.8333333333333
